$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I (I0) and J (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (style) of the existing header cell H1 onto I1:J1
# so they reuse the same bold/centered/bordered header style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Re-apply the values in case PasteSpecial touched them (it shouldn't,
# since we only paste formats, but this keeps things explicit/safe).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for rows 2-4 in columns I and J
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 9

$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 4
